$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 12: add acceptance info for "Pseudocode" task
$ws.Range("I12").Value = "Ablaufplan studieren"
$ws.Range("J12").Value = "akzeptiert"

# Row 13: add acceptance info for "Ablaufplan" task
$ws.Range("I13").Value = "Ablaufplan und Formeln erstellen"
$ws.Range("J13").Value = "akzeptiert"

# Row 16: status update + new "Endabnahme" task
$ws.Range("B16").Value = "in Arbeit"
$ws.Range("D16").Value = "Endabnahme"
$ws.Range("G16").Value = "alle"
$ws.Range("K16").Value = "-"
$ws.Range("L16").Value = "-"
$ws.Range("E16").Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("M16").Value = 40847

# Row 17: status update + new "Öffentlichkeitsarbeit" task
$ws.Range("B17").Value = "in Arbeit"
$ws.Range("D17").Value = "Öffentlichkeitsarbeit"

# Row 20: save-result task, extended description
$ws.Range("G20").Value = "Thomas, Mrosk/Schlufter"
$ws.Range("I20").Value = "Ergebnis speichern(xml oder txt)"
$ws.Range("K20").Value = "1h"

# sheet selection moved (cosmetic)
$ws.Range("D24").Select()
